$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 94 (duplicating the existing row 94), pushing rows
# 94-144 down to 95-145, so the sheet grows from A1:T144 to A1:T145.
$ws.Rows.Item(94).Copy()
$ws.Rows.Item(94).Insert()

# Overwrite the three fields that differ on the newly inserted record.
$ws.Range("D94").Value = 45001
$ws.Range("M94").Value = 100
$ws.Range("R94").Value = "Región del Maule"
